# Update countries & provincias Spain
#
# The sheet "Pais" lists COVID-19 stats per country, one row per country
# (rows 4-193ish), sorted by total cases descending. This refresh:
#   - bumps the "last updated" timestamp (row 1)
#   - updates total/new/active/dead counts for several countries whose
#     numbers changed (Suiza, Noruega, Canada, Estonia, Kazajistan, ...)
#   - inserts "Puerto Rico" as a new entry and re-sorts a block of
#     countries around it (rows 115-129), which shifts every row in that
#     block down by one position (new country name + that position's data)
#   - re-sorts another block around Aruba/Guyana/Surinam/Gabon (rows
#     144-150) the same way
#
# Rather than reproduce the sort algorithm, each affected row is written
# explicitly with its final (country name, B..H) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 21 de Marzo de 2020 a las 21:16' -> 'Datos actualizados a 21 de Marzo de 2020 a las 21:46'
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 21:46"

# Row 12: 'Suiza' -> 'Suiza'
$ws.Cells.Item(12, 2).Value = 6746
$ws.Cells.Item(12, 3).Value = 1131
$ws.Cells.Item(12, 5).Value = 6652
$ws.Cells.Item(12, 7).Value = 23
$ws.Cells.Item(12, 8).Value = 79

# Row 17: 'Noruega' -> 'Noruega'
$ws.Cells.Item(17, 2).Value = 2150
$ws.Cells.Item(17, 3).Value = 191
$ws.Cells.Item(17, 5).Value = 2137

# Row 21: 'Canada' -> 'Canada'
$ws.Cells.Item(21, 2).Value = 1280
$ws.Cells.Item(21, 3).Value = 193
$ws.Cells.Item(21, 5).Value = 1248

# Row 50: 'Estonia' -> 'Estonia'
$ws.Cells.Item(50, 6).Value = 0

# Row 94: 'Kazajistan' -> 'Kazajistan'
$ws.Cells.Item(94, 2).Value = 54
$ws.Cells.Item(94, 3).Value = 2
$ws.Cells.Item(94, 5).Value = 54

# Row 115: 'Cuba' -> 'Puerto Rico'
$ws.Cells.Item(115, 1).Value = "Puerto Rico"
$ws.Cells.Item(115, 3).Value = 7
$ws.Cells.Item(115, 5).Value = 21
$ws.Cells.Item(115, 8).Value = 0

# Row 116: 'Ghana' -> 'Cuba'
$ws.Cells.Item(116, 1).Value = "Cuba"
$ws.Cells.Item(116, 2).Value = 21
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 5).Value = 20
$ws.Cells.Item(116, 8).Value = 1

# Row 118: 'Jamaica' -> 'Ghana'
$ws.Cells.Item(118, 1).Value = "Ghana"
$ws.Cells.Item(118, 3).Value = 3
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(118, 5).Value = 18
$ws.Cells.Item(118, 7).Value = 1

# Row 119: 'Guayana Francesa' -> 'Jamaica'
$ws.Cells.Item(119, 1).Value = "Jamaica"
$ws.Cells.Item(119, 2).Value = 19
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 2
$ws.Cells.Item(119, 5).Value = 16
$ws.Cells.Item(119, 8).Value = 1

# Row 120: 'Monaco' -> 'Guayana Francesa'
$ws.Cells.Item(120, 1).Value = "Guayana Francesa"
$ws.Cells.Item(120, 3).Value = 3
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 18

# Row 121: 'Paraguay' -> 'Monaco'
$ws.Cells.Item(121, 1).Value = "Monaco"
$ws.Cells.Item(121, 3).Value = 7
$ws.Cells.Item(121, 4).Value = 1
$ws.Cells.Item(121, 6).Value = 0
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 0

# Row 122: 'Macao' -> 'Paraguay'
$ws.Cells.Item(122, 1).Value = "Paraguay"
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 17
$ws.Cells.Item(122, 6).Value = 1
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 1

# Row 123: 'Ruanda' -> 'Macao'
$ws.Cells.Item(123, 1).Value = "Macao"
$ws.Cells.Item(123, 2).Value = 18
$ws.Cells.Item(123, 3).Value = 1
$ws.Cells.Item(123, 4).Value = 10
$ws.Cells.Item(123, 5).Value = 8

# Row 124: 'Puerto Rico' -> 'Ruanda'
$ws.Cells.Item(124, 1).Value = "Ruanda"
$ws.Cells.Item(124, 3).Value = 0

# Row 127: 'Polinesia Francesa' -> 'Montenegro'
$ws.Cells.Item(127, 1).Value = "Montenegro"
$ws.Cells.Item(127, 2).Value = 16
$ws.Cells.Item(127, 3).Value = 2
$ws.Cells.Item(127, 5).Value = 16

# Row 128: 'Guam' -> 'Polinesia Francesa'
$ws.Cells.Item(128, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(128, 3).Value = 4

# Row 129: 'Montenegro' -> 'Guam'
$ws.Cells.Item(129, 1).Value = "Guam"
$ws.Cells.Item(129, 2).Value = 15
$ws.Cells.Item(129, 3).Value = 1
$ws.Cells.Item(129, 5).Value = 15

# Row 144: 'Aruba' -> 'Surinam'
$ws.Cells.Item(144, 1).Value = "Surinam"
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(144, 4).Value = 0
$ws.Cells.Item(144, 5).Value = 5

# Row 145: 'Guyana' -> 'Gabon'
$ws.Cells.Item(145, 1).Value = "Gabon"
$ws.Cells.Item(145, 3).Value = 1

# Row 146: 'San Martin (Parte Francesa)' -> 'Guyana'
$ws.Cells.Item(146, 1).Value = "Guyana"
$ws.Cells.Item(146, 2).Value = 5
$ws.Cells.Item(146, 8).Value = 1

# Row 147: 'Bahamas' -> 'Aruba'
$ws.Cells.Item(147, 1).Value = "Aruba"
$ws.Cells.Item(147, 2).Value = 5
$ws.Cells.Item(147, 4).Value = 1

# Row 148: 'Surinam' -> 'San Martin (Parte Francesa)'
$ws.Cells.Item(148, 1).Value = "San Martin (Parte Francesa)"

# Row 149: 'Nueva Caledonia' -> 'Bahamas'
$ws.Cells.Item(149, 1).Value = "Bahamas"
$ws.Cells.Item(149, 3).Value = 0

# Row 150: 'Gabon' -> 'Nueva Caledonia'
$ws.Cells.Item(150, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(150, 5).Value = 4
$ws.Cells.Item(150, 8).Value = 0
